$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text field updates (row 2) ---
$ws.Range("I2").Value = "Krish - Auto6"
$ws.Range("J2").Value = "Smith6"
$ws.Range("V2").Value = "John6"
$ws.Range("W2").Value = "Updik6"

# --- Checkbox selection moves from "race" (A2) to "family" (F2) ---
# Preserve the highlighted format by copying it along with the value.
$ws.Range("A2").Copy()
$ws.Range("F2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A2").Value = "No"
$ws.Range("F2").Value = "Yes"

$excel.CutCopyMode = 0

# --- Update the active selection state ---
$ws.Range("A2:XFD6").Select()
